$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.586.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.918.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.001'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4853'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.66%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2908'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06734'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '111.60'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.44'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.921.80'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07582'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.382'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6742'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.57%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '294.57'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.588.66'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.07'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007564'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.533'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.172.61'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.433'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.91%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.491'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '164.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.110'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1077'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.74%  '
$ws.Range("E30").Value = '  +3.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.143'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.096'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.35%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05025'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7415'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.139'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9997'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02032'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.696'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.84%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.691'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.023'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.94%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '109.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.35%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4454'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8669'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.891'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("E45").Value = '  +4.16%  '
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.274'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.06%  '
$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.39'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.63%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.245'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.36%  '
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1232'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.07%  '
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2553'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.00%  '
